# Update the "Förändrad" (Changed) date column (C) from 2023-09-06 (45175)
# to 2023-09-08 (45177) for all data rows (2 through 52) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newSerial = 45177

for ($row = 2; $row -le 52; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45175) {
        $cell.Value2 = $newSerial
    }
}
